$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 9,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.041827
$arr[0,3] = 0.125481
$arr[0,4] = 0.006279874897961605
$arr[0,5] = 0.006279874897961606
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.3252056666666667
$arr[0,9] = 0.975617
$arr[0,10] = 0.0158278498560244
$arr[0,11] = 0.0158278498560244
$arr[0,12] = 0.01360237741966667
$arr[0,13] = 0.122421396777
$arr[0,14] = 0.00009939691699955282
$arr[0,15] = 0.00009939691699955284
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 0.041827
$arr[1,3] = 0.125481
$arr[1,4] = 0.006279874897961605
$arr[1,5] = 0.006279874897961606
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 16.71131166666667
$arr[1,9] = 50.133935
$arr[1,10] = 0.8133441666880411
$arr[1,11] = 0.8133441666880411
$arr[1,12] = 0.6989840330816667
$arr[1,13] = 6.290856297735001
$arr[1,14] = 0.005107699615787728
$arr[1,15] = 0.005107699615787729
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 0.041827
$arr[2,3] = 0.125481
$arr[2,4] = 0.006279874897961605
$arr[2,5] = 0.006279874897961606
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 3.509903666666667
$arr[2,9] = 10.529711
$arr[2,10] = 0.1708279834559346
$arr[2,11] = 0.1708279834559346
$arr[2,12] = 0.1468087406656667
$arr[2,13] = 1.321278665991
$arr[2,14] = 0.001072778365174324
$arr[2,15] = 0.001072778365174324
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 1.127914
$arr[3,3] = 3.383742
$arr[3,4] = 0.1693441751896972
$arr[3,5] = 0.1693441751896972
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 0.3252056666666667
$arr[3,9] = 0.975617
$arr[3,10] = 0.0158278498560244
$arr[3,11] = 0.0158278498560244
$arr[3,12] = 0.3668040243126666
$arr[3,13] = 3.301236218814
$arr[3,14] = 0.002680354178894819
$arr[3,15] = 0.00268035417889482
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 1.127914
$arr[4,3] = 3.383742
$arr[4,4] = 0.1693441751896972
$arr[4,5] = 0.1693441751896972
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 16.71131166666667
$arr[4,9] = 50.133935
$arr[4,10] = 0.8133441666880411
$arr[4,11] = 0.8133441666880411
$arr[4,12] = 18.84892238719667
$arr[4,13] = 169.64030148477
$arr[4,14] = 0.1377350970531379
$arr[4,15] = 0.1377350970531379
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 1.127914
$arr[5,3] = 3.383742
$arr[5,4] = 0.1693441751896972
$arr[5,5] = 0.1693441751896972
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 3.509903666666667
$arr[5,9] = 10.529711
$arr[5,10] = 0.1708279834559346
$arr[5,11] = 0.1708279834559346
$arr[5,12] = 3.958869484284666
$arr[5,13] = 35.629825358562
$arr[5,14] = 0.02892872395766448
$arr[5,15] = 0.02892872395766448
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 5.490742
$arr[6,3] = 16.472226
$arr[6,4] = 0.8243759499123412
$arr[6,5] = 0.8243759499123412
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 0.3252056666666667
$arr[6,9] = 0.975617
$arr[6,10] = 0.0158278498560244
$arr[6,11] = 0.0158278498560244
$arr[6,12] = 1.785620412604667
$arr[6,13] = 16.070583713442
$arr[6,14] = 0.01304809876013003
$arr[6,15] = 0.01304809876013003
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 5.490742
$arr[7,3] = 16.472226
$arr[7,4] = 0.8243759499123412
$arr[7,5] = 0.8243759499123412
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 16.71131166666667
$arr[7,9] = 50.133935
$arr[7,10] = 0.8133441666880411
$arr[7,11] = 0.8133441666880411
$arr[7,12] = 91.75750084325666
$arr[7,13] = 825.81750758931
$arr[7,14] = 0.6705013700191155
$arr[7,15] = 0.6705013700191155
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 5.490742
$arr[8,3] = 16.472226
$arr[8,4] = 0.8243759499123412
$arr[8,5] = 0.8243759499123412
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 3.509903666666667
$arr[8,9] = 10.529711
$arr[8,10] = 0.1708279834559346
$arr[8,11] = 0.1708279834559346
$arr[8,12] = 19.27197547852067
$arr[8,13] = 173.447779306686
$arr[8,14] = 0.1408264811330958
$arr[8,15] = 0.1408264811330958

$range = $ws.Range("E2:T10")
$range.Value = $arr

Write-Output "Done updating Reln-Vldlr sheet"